# Credit simulation data update:
# - new assertion columns (assert_request_amount / assert_rut)
# - refreshed test user / amount / dues / rut sample data
# - new fonts/styles for the ci_document + assert_request_amount cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the two extra assertion columns
$ws.Range("O1").Value = "assert_request_amount"
$ws.Range("P1").Value = "assert_rut"

# Refreshed row 2 sample data
$ws.Range("A2").Value = "55589143"
$ws.Range("B2").Value = "QA2022"
$ws.Range("C2").Value = "1500000"
$ws.Range("D2").Value = "48"
$ws.Range("E2").Value = "Sin meses de gracia"
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = "$ 1.500.000 - 48 Cuotas"
$ws.Range("H2").Value = 103238224
$ws.Range("I2").Value = "PROVIDENCIA 123"
$ws.Range("J2").Value = "Bio-Bio"
$ws.Range("K2").Value = "Arauco"
$ws.Range("L2").Value = "Itau"
$ws.Range("M2").Value = "Cuenta Ahorro"
$ws.Range("N2").Value = 12345678
$ws.Range("O2").Value = "$ 1.500.000"
$ws.Range("P2").Value = "5.558.914-3"

# Distinguish the id-document style and the new assertion-amount style with
# their own font colors (dark black / dark gray)
$ws.Range("H2").Font.Color = 0
$ws.Range("O2").Font.Color = 4473924

# New columns need to fit their (wider) content
$ws.Columns.Item(15).AutoFit() | Out-Null
$ws.Columns.Item(16).AutoFit() | Out-Null

# Move the active selection like a user would after filling this row in
$ws.Range("H9").Select() | Out-Null
